$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Rows.Item(245).Insert()

$ws.Range("A245").Value = 11
$ws.Range("B245").Value = "Vega Monumental Concepción"
$ws.Range("C245").Value = "Bíobío"
$ws.Range("D245").Value = 44504
$ws.Range("E245").Value = 8
$ws.Range("F245").Value = 100112020
$ws.Range("G245").Value = "Tomate"
$ws.Range("H245").Value = "Larga vida"
$ws.Range("I245").Value = "Primera"
$ws.Range("J245").Value = 1250
$ws.Range("K245").Value = 7500
$ws.Range("L245").Value = 8000
$ws.Range("M245").Value = 7760
$ws.Range("N245").Value = "`$/caja 12 kilos"
$ws.Range("O245").Value = "Región de Arica y Parinacota"
$ws.Range("P245").Value = 647
$ws.Range("Q245").Value = 12
$ws.Range("R245").Value = "Hortaliza"
